# Updated filter with line breaks on methods
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

# Re-wrap the "C." (85%, 65%, 50% ...) and "E." (Low Lake Powell ...) method
# descriptions so the line breaks fall in different places. ":" is used in
# this workbook as the in-cell line-break marker, matching the existing
# text in the other rows.
$ws.Range("A4").Value = "C. 85%, 65%, and 50% of:2000 to 2018 average:flow (2022)"
$ws.Range("A6").Value = "E. Low Lake Powell:releases + gains through:Grand Canyon (2025)"

# Move the active selection from A8 to A7.
$ws.Range("A7").Select()
